# Auto-generated Excel COM-interop script
# Applies scheduled-runner market-data refresh to the 8 Leve-profit sheets
# (columns H:N = currentAveragePrice*, LevePriceNQ/HQ, LeveProfitNQ/HQ)
$wb = $excel.ActiveWorkbook

# --- Sheet 1: ALC ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("H18").Value = 111111310
$ws.Range("I18").Value = 111111310
$ws.Range("K18").Value = 111111310
$ws.Range("M18").Value = -111111026
$ws.Range("H41").Value = 1436.5714
$ws.Range("J41").Value = 2259
$ws.Range("L41").Value = 2259
$ws.Range("N41").Value = -3139
$ws.Range("H51").Value = 5387.7144
$ws.Range("J51").Value = 6893.5557
$ws.Range("L51").Value = 6893.5557
$ws.Range("N51").Value = -7861.5557
$ws.Range("H92").Value = 376.70587
$ws.Range("I92").Value = 281.1875
$ws.Range("K92").Value = 281.1875
$ws.Range("M92").Value = 966.8125
$ws.Range("H116").Value = 21319454
$ws.Range("I116").Value = 35791052
$ws.Range("J116").Value = 14290391
$ws.Range("K116").Value = 35791052
$ws.Range("L116").Value = 14290391
$ws.Range("M116").Value = -35787610
$ws.Range("N116").Value = -14297275
$ws.Range("H132").Value = 8583.387000000001
$ws.Range("I132").Value = 1827.3334
$ws.Range("K132").Value = 5482.0002
$ws.Range("M132").Value = -2952.0002
$ws.Range("H137").Value = 7940002.5
$ws.Range("I137").Value = 1433.7727
$ws.Range("J137").Value = 16672428
$ws.Range("K137").Value = 4301.3181
$ws.Range("L137").Value = 50017284
$ws.Range("M137").Value = -1751.3181
$ws.Range("N137").Value = -50022384
$ws.Range("H138").Value = 7009
$ws.Range("J138").Value = 7749.1143
$ws.Range("L138").Value = 23247.3429
$ws.Range("N138").Value = -33527.3429
$ws.Range("H141").Value = 6129.76
$ws.Range("I141").Value = 5358.7827
$ws.Range("K141").Value = 16076.3481
$ws.Range("M141").Value = -10896.3481

# --- Sheet 2: ARM ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("H4").Value = 300
$ws.Range("I4").Value = 300
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 300
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -184
$ws.Range("N4").ClearContents()
$ws.Range("H5").Value = 710
$ws.Range("I5").Value = 710
$ws.Range("K5").Value = 710
$ws.Range("M5").Value = -598
$ws.Range("H32").Value = 3873.8064
$ws.Range("I32").Value = 2240.3208
$ws.Range("K32").Value = 2240.3208
$ws.Range("M32").Value = -1953.3208
$ws.Range("H74").Value = 15626558
$ws.Range("I74").Value = 25001318
$ws.Range("K74").Value = 25001318
$ws.Range("M74").Value = -25000444
$ws.Range("H77").Value = 15626558
$ws.Range("I77").Value = 25001318
$ws.Range("K77").Value = 125006590
$ws.Range("M77").Value = -125002222

# --- Sheet 3: BSM ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("H4").Value = 710
$ws.Range("I4").Value = 710
$ws.Range("K4").Value = 710
$ws.Range("M4").Value = -595
$ws.Range("H22").Value = 440.58334
$ws.Range("I22").Value = 440.58334
$ws.Range("K22").Value = 440.58334
$ws.Range("M22").Value = -267.58334
$ws.Range("H80").Value = 326.36365
$ws.Range("J80").Value = 234.85715
$ws.Range("L80").Value = 234.85715
$ws.Range("N80").Value = -2230.85715
$ws.Range("H83").Value = 326.36365
$ws.Range("J83").Value = 234.85715
$ws.Range("L83").Value = 1174.28575
$ws.Range("N83").Value = -11158.28575
$ws.Range("H86").Value = 1521.875
$ws.Range("I86").Value = 1391
$ws.Range("K86").Value = 1391
$ws.Range("M86").Value = -268
$ws.Range("H89").Value = 1521.875
$ws.Range("I89").Value = 1391
$ws.Range("K89").Value = 6955
$ws.Range("M89").Value = -1339
$ws.Range("H107").Value = 1339.6111
$ws.Range("I107").Value = 1533.6666
$ws.Range("K107").Value = 1533.6666
$ws.Range("M107").Value = 386.3334

# --- Sheet 4: CRP ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("H7").Value = 244.4
$ws.Range("I7").Value = 162.33333
$ws.Range("K7").Value = 162.33333
$ws.Range("M7").Value = -49.33332999999999
$ws.Range("H31").Value = 2078
$ws.Range("I31").Value = 1266.1666
$ws.Range("J31").Value = 4165.5713
$ws.Range("K31").Value = 1266.1666
$ws.Range("L31").Value = 4165.5713
$ws.Range("M31").Value = -971.1666
$ws.Range("N31").Value = -4755.5713
$ws.Range("H34").Value = 2078
$ws.Range("I34").Value = 1266.1666
$ws.Range("J34").Value = 4165.5713
$ws.Range("K34").Value = 1266.1666
$ws.Range("L34").Value = 4165.5713
$ws.Range("M34").Value = -1064.1666
$ws.Range("N34").Value = -4569.5713
$ws.Range("H105").Value = 1058.1786
$ws.Range("I105").Value = 958
$ws.Range("K105").Value = 958
$ws.Range("M105").Value = 789

# --- Sheet 5: CUL ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("H5").Value = 713.6
$ws.Range("J5").Value = 600
$ws.Range("L5").Value = 1800
$ws.Range("N5").Value = -2024
$ws.Range("H86").Value = 212.4
$ws.Range("I86").Value = 191
$ws.Range("J86").Value = 226.66667
$ws.Range("K86").Value = 573
$ws.Range("L86").Value = 680.00001
$ws.Range("M86").Value = 613
$ws.Range("N86").Value = -3052.00001
$ws.Range("H89").Value = 212.4
$ws.Range("I89").Value = 191
$ws.Range("J89").Value = 226.66667
$ws.Range("K89").Value = 1719
$ws.Range("L89").Value = 2040.00003
$ws.Range("M89").Value = 4209
$ws.Range("N89").Value = -13896.00003
$ws.Range("H122").Value = 856.2
$ws.Range("J122").Value = 795.5
$ws.Range("L122").Value = 7159.5
$ws.Range("N122").Value = -12059.5
$ws.Range("H126").Value = 12626.5
$ws.Range("I126").Value = 6439.75
$ws.Range("K126").Value = 19319.25
$ws.Range("M126").Value = -14379.25
$ws.Range("H129").Value = 1997.2222
$ws.Range("J129").Value = 2136.6
$ws.Range("L129").Value = 6409.799999999999
$ws.Range("N129").Value = -16409.8
$ws.Range("H135").Value = 713.6
$ws.Range("J135").Value = 600
$ws.Range("L135").Value = 5400
$ws.Range("N135").Value = -10470

# --- Sheet 6: GSM ---
$ws = $wb.Worksheets.Item(6)
$ws.Range("H2").Value = 247.8125
$ws.Range("I2").Value = 310
$ws.Range("J2").Value = 199.44444
$ws.Range("K2").Value = 310
$ws.Range("L2").Value = 199.44444
$ws.Range("M2").Value = -197
$ws.Range("N2").Value = -425.44444
$ws.Range("H80").Value = 15613.286
$ws.Range("I80").Value = 5000
$ws.Range("K80").Value = 5000
$ws.Range("M80").Value = -4002
$ws.Range("H83").Value = 15613.286
$ws.Range("I83").Value = 5000
$ws.Range("K83").Value = 25000
$ws.Range("M83").Value = -20008
$ws.Range("H134").Value = 899505.3
$ws.Range("J134").Value = 899505.3
$ws.Range("L134").Value = 2698515.9
$ws.Range("N134").Value = -2703585.9

# --- Sheet 7: LTW ---
$ws = $wb.Worksheets.Item(7)
$ws.Range("H18").Value = 14999
$ws.Range("J18").Value = 14999
$ws.Range("L18").Value = 14999
$ws.Range("N18").Value = -15343
$ws.Range("H61").Value = 3000
$ws.Range("J61").Value = 3000
$ws.Range("L61").Value = 3000
$ws.Range("N61").Value = -3404
$ws.Range("H68").Value = 1751812.5
$ws.Range("I68").Value = 2275416.8
$ws.Range("K68").Value = 2275416.8
$ws.Range("M68").Value = -2274667.8
$ws.Range("H71").Value = 1751812.5
$ws.Range("I71").Value = 2275416.8
$ws.Range("K71").Value = 11377084
$ws.Range("M71").Value = -11373340
$ws.Range("H113").Value = 3000
$ws.Range("J113").Value = 3000
$ws.Range("L113").Value = 3000
$ws.Range("N113").Value = -7340
$ws.Range("H134").Value = 20000
$ws.Range("J134").Value = 20000
$ws.Range("L134").Value = 20000
$ws.Range("N134").Value = -30140

# --- Sheet 8: WVR ---
$ws = $wb.Worksheets.Item(8)
$ws.Range("H100").Value = 953958.7
$ws.Range("I100").Value = 1429519
$ws.Range("J100").Value = 2838.1428
$ws.Range("K100").Value = 2859038
$ws.Range("L100").Value = 5676.2856
$ws.Range("M100").Value = -2858497
$ws.Range("N100").Value = -6758.2856
$ws.Range("H107").Value = 5169.2
$ws.Range("J107").Value = 3301.5
$ws.Range("L107").Value = 9904.5
$ws.Range("N107").Value = -13744.5
$ws.Range("H132").Value = 10686004
$ws.Range("I132").Value = 1588986
$ws.Range("J132").Value = 29415158
$ws.Range("K132").Value = 4766958
$ws.Range("L132").Value = 88245474
$ws.Range("M132").Value = -4764428
$ws.Range("N132").Value = -88250534
$ws.Range("H136").Value = 7666.6123
$ws.Range("I136").Value = 3086.2964
$ws.Range("J136").Value = 9408.423000000001
$ws.Range("K136").Value = 9258.889200000001
$ws.Range("L136").Value = 28225.269
$ws.Range("M136").Value = -6708.889200000001
$ws.Range("N136").Value = -33325.269
